# Auto-generated edit script applying the Marilith_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H,I,J,K,L,M,N)
# for a set of leve rows across multiple crafting-job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2509.0476
$ws.Range("J17").Value = 2509.0476
$ws.Range("L17").Value = 7527.1428
$ws.Range("N17").Value = -7863.1428

$ws.Range("H69").Value = 3000
$ws.Range("J69").Value = 3000
$ws.Range("L69").Value = 9000
$ws.Range("N69").Value = -10748

$ws.Range("H72").Value = 3000
$ws.Range("J72").Value = 3000
$ws.Range("L72").Value = 27000
$ws.Range("N72").Value = -35736

$ws.Range("H76").Value = 4066.6
$ws.Range("J76").Value = 4666.6665
$ws.Range("L76").Value = 4666.6665
$ws.Range("N76").Value = -5296.6665

$ws.Range("H79").Value = 4066.6
$ws.Range("J79").Value = 4666.6665
$ws.Range("L79").Value = 4666.6665
$ws.Range("N79").Value = -6850.6665

$ws.Range("H112").Value = 3383.1667
$ws.Range("J112").Value = 7500
$ws.Range("L112").Value = 22500
$ws.Range("N112").Value = -24716

$ws.Range("H138").Value = 2197
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12081.8
$ws.Range("I32").Value = 12946.667
$ws.Range("K32").Value = 12946.667
$ws.Range("M32").Value = -12659.667

$ws.Range("H61").Value = 3665.8333
$ws.Range("J61").Value = 3999.5
$ws.Range("L61").Value = 3999.5
$ws.Range("N61").Value = -4423.5

$ws.Range("H88").Value = 2881.3076
$ws.Range("I88").Value = 914.4
$ws.Range("J88").Value = 4110.625
$ws.Range("K88").Value = 914.4
$ws.Range("L88").Value = 4110.625
$ws.Range("M88").Value = -508.4
$ws.Range("N88").Value = -4922.625

$ws.Range("H91").Value = 2881.3076
$ws.Range("I91").Value = 914.4
$ws.Range("J91").Value = 4110.625
$ws.Range("K91").Value = 914.4
$ws.Range("L91").Value = 4110.625
$ws.Range("M91").Value = 489.6
$ws.Range("N91").Value = -6918.625

$ws.Range("H122").Value = 1200
$ws.Range("I122").Value = 1200
$ws.Range("K122").Value = 3600
$ws.Range("M122").Value = -1150

$ws.Range("H136").Value = 3665.8333
$ws.Range("J136").Value = 3999.5
$ws.Range("L136").Value = 11998.5
$ws.Range("N136").Value = -17098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5439.0625
$ws.Range("I20").Value = 4935.5454
$ws.Range("K20").Value = 4935.5454
$ws.Range("M20").Value = -4688.5454

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws.Range("H107").Value = 843.5
$ws.Range("I107").Value = 715
$ws.Range("K107").Value = 715
$ws.Range("M107").Value = 1205

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2539.0244
$ws.Range("I31").Value = 2066.1428
$ws.Range("K31").Value = 2066.1428
$ws.Range("M31").Value = -1771.1428

$ws.Range("H34").Value = 2539.0244
$ws.Range("I34").Value = 2066.1428
$ws.Range("K34").Value = 2066.1428
$ws.Range("M34").Value = -1864.1428

$ws.Range("H58").Value = 2855
$ws.Range("I58").Value = 3060
$ws.Range("J58").Value = 2547.5
$ws.Range("K58").Value = 3060
$ws.Range("L58").Value = 2547.5
$ws.Range("M58").Value = -2857
$ws.Range("N58").Value = -2953.5

$ws.Range("H62").Value = 4598.8
$ws.Range("I62").Value = 4498
$ws.Range("J62").Value = 4750
$ws.Range("K62").Value = 4498
$ws.Range("L62").Value = 4750
$ws.Range("M62").Value = -3874
$ws.Range("N62").Value = -5998

$ws.Range("H65").Value = 4598.8
$ws.Range("I65").Value = 4498
$ws.Range("J65").Value = 4750
$ws.Range("K65").Value = 22490
$ws.Range("L65").Value = 23750
$ws.Range("M65").Value = -19370
$ws.Range("N65").Value = -29990

$ws.Range("H99").Value = 4874.75
$ws.Range("I99").Value = 4863.4546
$ws.Range("K99").Value = 4863.4546
$ws.Range("M99").Value = -3365.4546

$ws.Range("H126").Value = 4874.75
$ws.Range("I126").Value = 4863.4546
$ws.Range("K126").Value = 14590.3638
$ws.Range("M126").Value = -12120.3638

$ws.Range("H136").Value = 2855
$ws.Range("I136").Value = 3060
$ws.Range("J136").Value = 2547.5
$ws.Range("K136").Value = 9180
$ws.Range("L136").Value = 7642.5
$ws.Range("M136").Value = -6630
$ws.Range("N136").Value = -12742.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2908.75
$ws.Range("J5").Value = 2095.4
$ws.Range("L5").Value = 6286.200000000001
$ws.Range("N5").Value = -6510.200000000001

$ws.Range("H34").Value = 898.0769
$ws.Range("I34").Value = 784.5
$ws.Range("J34").Value = 1079.8
$ws.Range("K34").Value = 2353.5
$ws.Range("L34").Value = 3239.4
$ws.Range("M34").Value = -2269.5
$ws.Range("N34").Value = -3407.4

$ws.Range("H113").Value = 1197.5
$ws.Range("J113").Value = 1363.7778
$ws.Range("L113").Value = 4091.3334
$ws.Range("N113").Value = -8431.3334

$ws.Range("H135").Value = 2908.75
$ws.Range("J135").Value = 2095.4
$ws.Range("L135").Value = 18858.6
$ws.Range("N135").Value = -23928.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H68").Value = 29662.5
$ws.Range("I68").Value = 2883.3333
$ws.Range("K68").Value = 2883.3333
$ws.Range("M68").Value = -2134.3333

$ws.Range("H71").Value = 29662.5
$ws.Range("I71").Value = 2883.3333
$ws.Range("K71").Value = 14416.6665
$ws.Range("M71").Value = -10672.6665

$ws.Range("H94").Value = 57000
$ws.Range("I94").Value = 3000
$ws.Range("J94").Value = 75000
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 75000
$ws.Range("M94").Value = -2324
$ws.Range("N94").Value = -76352

$ws.Range("H122").Value = 2851.6
$ws.Range("I122").Value = 3203.2
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 9609.599999999999
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -7159.599999999999
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1698
$ws.Range("I122").Value = 1698
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5094
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2644
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 1560.2727
$ws.Range("I132").Value = 1406.5
$ws.Range("K132").Value = 4219.5
$ws.Range("M132").Value = -1689.5
